$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns, plus swapped rows' Coin/Link (B/C)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.031.99'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.98%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.720.01'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.34%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9954'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3776'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '50.07'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3502'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.197'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07482'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9989'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.339'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.98'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.003'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.720.77'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001126'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9966'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06680'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '84.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.36'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.399'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.40'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +11.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '25.099.25'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.431'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.821'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.55'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '150.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '132.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.911.12'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.60%  '
$ws.Range('E31').Value = '  +20.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.884'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.209'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.16%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '13.81'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.77%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08848'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.773'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.660'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.48%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06586'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.81%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02430'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.977'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2213'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.15%  '
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6465'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9965'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.01'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6161'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.828'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.150'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '129.60'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.99%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07330'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '80.11'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.50%  '
